# "viz y tablas update"
# Update the "Ficha técnica" sheet: remove the DIMENSIÓN/Accesibilidad row
# (shifting CONINDICADOR/NOMINDICADOR/DEFINICIÓN/CÁLCULO up by one row) and
# append two new metadata rows: TIPOIND/Resultados and CITA/<source note>.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

$ws.Range("A3").Value = "CONINDICADOR"
$ws.Range("B3").Value = "Porcentaje de ocupados que no aporta a la seguridad social"

$ws.Range("A4").Value = "NOMINDICADOR"
$ws.Range("B4").Value = "Porcentaje de ocupados que no aporta a la seguridad social"

$ws.Range("A5").Value = "DEFINICIÓN"
$ws.Range("B5").Value = "El indicador mide el porcentaje de ocupados que no aporta a la seguridad social."

$ws.Range("A6").Value = "CÁLCULO"
$ws.Range("B6").Value = "Para cada año calcular: (Cantidad de ocupados que no aporta a la seguridad social / Cantidad de ocupados)*100"

$ws.Range("A7").Value = "TIPOIND"
$ws.Range("B7").Value = "Resultados"

$ws.Range("A8").Value = "CITA"
$ws.Range("B8").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"
